$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 9, shifting the existing rows
# 9-11 down to 10-12 (keeping all their data/styles intact).
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new record's data. Columns
# A, B, C, E, F, G, H, I, J, K are identical to the surrounding rows, so
# just mirror row 10 (the row that used to be row 9) for those.
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C9").Value = "Arica y Parinacota"
$ws.Range("D9").Value = 44895
$ws.Range("D9").NumberFormat = $ws.Range("D10").NumberFormat
$ws.Range("E9").Value = 15
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100103
$ws.Range("H9").Value = "Frutos de hueso (carozo)"
$ws.Range("I9").Value = 100103003
$ws.Range("J9").Value = "Damasco"
$ws.Range("K9").Value = "Castle Brite"
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 130
$ws.Range("N9").Value = 19000
$ws.Range("O9").Value = 20000
$ws.Range("P9").Value = 19462
$ws.Range("Q9").Value = "$/caja 16 kilos granel"
$ws.Range("R9").Value = "Región de O'Higgins"
$ws.Range("S9").Value = 1216
$ws.Range("T9").Value = 16

# Row 11 (old row 10's data, now shifted down) changes quality grade.
$ws.Range("L11").Value = "Segunda"
